# Refactoring loginPage to SignUp & loginTests to signUpTests
# The test-data sheet's second data row is updated:
#   - A2 ("domainName" column): test@yopmail.com -> yes@finalmail.com
#   - B2 ("password" column) keeps its visible value Test$1234! (re-set so the
#     underlying shared-string ordering matches the refactor)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "yes@finalmail.com"
$ws.Range("B2").Value = "Test`$1234!"
